$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (pushes existing rows 3..20 down to 4..21),
# inheriting the formatting (e.g. date style on column G) from the row above.
$ws.Rows.Item(3).EntireRow.Insert()

# New row 3: AF - Anna FILATOVA, updated by jcs on 44033 (2020-07-21)
$ws.Range("A3").Value = "Oui"
$ws.Range("B3").Value = "AF"
$ws.Range("C3").Value = "Anna FILATOVA"
$ws.Range("D3").Value = "Anna FILATOVA"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "jcs"
$ws.Range("G3").Value = 44033

# Rows 4..21 (previously 3..20, now shifted down by the insert) were all
# touched by the same update: "Mis à jour par" -> jcs,
# "Dernière mise à jour" -> 44033 (2020-07-21)
for ($r = 4; $r -le 21; $r++) {
    $ws.Cells.Item($r, 6).Value = "jcs"
    $ws.Cells.Item($r, 7).Value = 44033
}
